# Auto-generated edit script: update cryptos list data per commit diff
# (Sat Jul 27 10:16:52 UTC 2024 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose refreshed value is a plain decimal number must be
# forced to Text so Excel doesn't silently coerce them to numeric cells
# (these prices are stored as text throughout the sheet, e.g. '68.253.29').
$textCells = @("D5", "D6", "D10", "D14", "D16", "D19", "D20", "D21", "D23", "D26", "D27", "D29", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49")
foreach ($r in $textCells) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.253.29'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '3.278.63'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '588.84'
$ws.Range('E5').Value = '  +1.95%  '
$ws.Range('D6').Value = '186.30'
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  +4.75%  '
$ws.Range('D10').Value = '6.74'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').Value = '3.846.03'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '28.64'
$ws.Range('E14').Value = '  +2.12%  '
$ws.Range('D15').Value = '68.246.59'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '0.0000172'
$ws.Range('D17').Value = '3.273.64'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '13.67'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').Value = '382.31'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('D21').Value = '7.75'
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = '71.53'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +2.58%  '
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').Value = '0.191'
$ws.Range('E26').Value = '  +6.58%  '
$ws.Range('D27').Value = '9.79'
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').Value = '5.84'
$ws.Range('E29').Value = '  +3.65%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('E32').Value = '  +4.68%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +3.19%  '
$ws.Range('D36').Value = '163.02'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').Value = '6.80'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '2.69'
$ws.Range('E40').Value = '  +2.60%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '26.74'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '4.61'
$ws.Range('E42').Value = '  +4.71%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '0.0694'
$ws.Range('E43').Value = '  +2.94%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '25.58'
$ws.Range('E44').Value = '  -1.41%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '41.40'
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '345.01'
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.642.47'
$ws.Range('E47').Value = '  -4.59%  '
$ws.Range('D48').Value = '0.0285'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('D49').Value = '32.46'
$ws.Range('E49').Value = '  +4.69%  '
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('E51').Value = '  -0.12%  '

Write-Host "Applied 91 cell updates across the cryptos list."
